# Add messy header challenge data
# - Fix the messy "Species Name " header on the Data sheet to "Species Name °C"
# - Clear the stray average-weight value for "Apis mellifera" (C6)
# - Update view state (selections / zoom) to match the saved workbook

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsData = $wb.Worksheets.Item("Data")

# -- Data sheet edits --------------------------------------------------
$wsData.Activate()
$wsData.Range("A4").Value = "Species Name °C"
$wsData.Range("C6").ClearContents()

# -- View / selection state --------------------------------------------
# Metadata sheet keeps a selection over A1:D3 (no data there beyond A1:A3)
$null = $wsMeta.Range("A1:D3").Select()

# Data sheet stays the active tab, zoomed to 205% with D8 selected
$wsData.Activate()
$excel.ActiveWindow.Zoom = 205
$null = $wsData.Range("D8").Select()
